$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows for "RM 232" (row 26) and "SC 92" (row 28).
# Delete the higher-numbered row first so the other row index is unaffected.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Scalar cell value changes (newly filled-in / newly blanked values).
$ws.Range("D3").Value = -14.2
$ws.Range("F4").Value = $null
$ws.Range("D5").Value = $null
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F13").Value = $null
$ws.Range("F14").Value = $null
$ws.Range("D21").Value = -14.3
$ws.Range("D23").Value = $null

# "SC 193" row (now row 32) had its D value filled in.
$ws.Range("D32").Value = -14.7
